$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate {
    param($Range, $DateText)
    # Cells holding "Latest Date" values are plain text (e.g. "2025-09-01"),
    # not real Excel dates. Writing an ISO-looking string straight into a
    # General-formatted cell would cause Excel to auto-convert it into a
    # date serial number, so we force Text format first, assign the literal
    # string, then restore a General display format (the underlying value
    # stays text either way).
    $Range.NumberFormat = "@"
    $Range.Value = $DateText
    $Range.NumberFormat = "General"
}

# --- Row 28: UMCSENT (Mich NTM Inflation Exp) ---
Set-TextDate $ws.Range("N28") "2025-09-01"
$ws.Range("Q28").Value = 55.1
$ws.Range("R28").Value = 58.2
$ws.Range("S28").Value = 61.7
$ws.Range("T28").Value = 60.7
$ws.Range("U28").Value = 52.2

# --- Row 29: T5YIFR (5yr, 5yr Forward) ---
Set-TextDate $ws.Range("N29") "2025-10-24"

# --- Row 30: T10YIE (10yr TIPS) ---
Set-TextDate $ws.Range("N30") "2025-10-24"
$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.3
$ws.Range("S30").Value = 2.29
$ws.Range("T30").Value = 2.28
$ws.Range("U30").Value = 2.26

# --- Row 48: DGS2 (2y UST) ---
Set-TextDate $ws.Range("N48") "2025-10-23"
$ws.Range("Q48").Value = 3.48
$ws.Range("R48").Value = 3.45
$ws.Range("S48").Value = 3.45
$ws.Range("T48").Value = 3.46
$ws.Range("U48").Value = 3.46

# --- Row 49: DGS5 (5y UST) ---
Set-TextDate $ws.Range("N49") "2025-10-23"
$ws.Range("Q49").Value = 3.61
$ws.Range("R49").Value = 3.56
$ws.Range("S49").Value = 3.56
$ws.Range("T49").Value = 3.58
$ws.Range("U49").Value = 3.59

# --- Row 50: DGS10 (10y UST) ---
Set-TextDate $ws.Range("N50") "2025-10-23"
$ws.Range("Q50").Value = 4.01
$ws.Range("R50").Value = 3.97
$ws.Range("S50").Value = 3.98
$ws.Range("T50").Value = 4
$ws.Range("U50").Value = 4.02

# --- Row 52: DBAA (BAA) ---
Set-TextDate $ws.Range("N52") "2025-10-23"
$ws.Range("Q52").Value = 5.67
$ws.Range("R52").Value = 5.66
$ws.Range("S52").Value = 5.65
$ws.Range("T52").Value = 5.68
$ws.Range("U52").Value = 5.72
